$d = $word.ActiveDocument

function SplitRange($doc, $startPos, $endPos) {
    # Forces a run boundary at startPos/endPos without altering the visible
    # text or its effective formatting: toggling Bold on then off leaves the
    # run split in place (Word/this engine won't silently re-merge runs that
    # were independently touched).
    $rr = $doc.Range($startPos, $endPos)
    $rr.Font.Bold = 1
    $rr.Font.Bold = 0
}

# ---------------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark from its original spot (between
#    "open-a" and "ccess"); it will be re-created later at its new location.
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# ---------------------------------------------------------------------------
# 2) Introduction paragraph: "... at the end of studies CTUs. In particular"
#    becomes "... at the end of studies on clinical trials units. In particular"
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$oldPhrase = "load at the end of studies CTUs. In particular"
$idx = $full.IndexOf($oldPhrase)
$len = "load at the end of studies CTUs".Length
$r = $d.Range($idx, $idx + $len)
$r.Text = "load at the end of studies on clinical trials units"

# Re-split the resulting (merged) run back into the five pieces shown by the
# canonical diff: "ver, imposes a large extra work" | "load at the end of
# studies" | " on clinical trials units" | ". In particular, the
# adverse-event-" | "reporting component requires entering:"
$full = $d.Content.Text
$workIdx = $full.IndexOf("ver, imposes a large extra work")
$b1 = $workIdx + "ver, imposes a large extra work".Length
$b2 = $b1 + "load at the end of studies".Length
$b3 = $b2 + " on clinical trials units".Length
$b4 = $b3 + ". In particular, the adverse-event-".Length
$b5 = $b4 + "reporting component requires entering:".Length

SplitRange $d $b4 $b5
SplitRange $d $b3 $b4
SplitRange $d $b2 $b3
SplitRange $d $b1 $b2

# ---------------------------------------------------------------------------
# 3) "Potential Relevance and Impact" paragraph: "... manually entering
#    1000s of data points using ..." becomes "... manually entering a large
#    amount of data points (e.g. over 1000 datum points for a recent
#    oncology study) using ..."
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$oldPhrase2 = "manually entering 1000s of data points using"
$newPhrase2 = "manually entering a large amount of data points (e.g. over 1000 datum points for a recent oncology study) using"
$idx2 = $full.IndexOf($oldPhrase2)
$r2 = $d.Range($idx2, $idx2 + $oldPhrase2.Length)
$r2.Text = $newPhrase2

# Re-split into the four pieces shown by the canonical diff: "The tool
# should remove the workload on CTUs of manually entering " | "a large
# amount of " | "data points" | " (e.g. over 1000 datum points for a recent
# oncology study)"
$full = $d.Content.Text
$pStart = $full.IndexOf("The tool should remove the workload on CTUs")
$c1 = $pStart + "The tool should remove the workload on CTUs of manually entering ".Length
$c2 = $c1 + "a large amount of ".Length
$c3 = $c2 + "data points".Length
$c4 = $c3 + " (e.g. over 1000 datum points for a recent oncology study)".Length

SplitRange $d $c3 $c4
SplitRange $d $c2 $c3
SplitRange $d $c1 $c2

# ---------------------------------------------------------------------------
# 4) Re-create the "_GoBack" bookmark right after the new
#    "... oncology study)" text, immediately before " using the
#    web-interface".
# ---------------------------------------------------------------------------
$rBookmark = $d.Range($c4, $c4)
[void]$d.Bookmarks.Add("_GoBack", $rBookmark)

Write-Output "Edit complete."
